$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 4 (shifts rows 4:23 down to 5:24)
$ws.Rows.Item(4).Insert()

# Fill in the new row 4 values: px = 3, em = $C$2*B4 (not part of the shared formula group)
$ws.Range("B4").Value = 3
$ws.Range("C4").Formula = "=`$C`$2*B4"

# Update sheet view: scroll back to top-left, select the new inserted row's C cells
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("C3:C4").Select()
